# Update server results across the yearly scenario sheets
$wb = $excel.ActiveWorkbook

# --- 2025 ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 132.021302159999
$ws.Range("E2").Value = 29438.60104869408
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 14174.76361456303
$ws.Range("L2").Value = 52558.32829870572
$ws.Range("M2").Value = 11132.49242710001
$ws.Range("N2").Value = 7211.629754433677
$ws.Range("O2").Value = 6821.594717536615

# --- 2030 ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 4743.52105312454
$ws.Range("E2").Value = 56948.83413643156
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 32800.48841919314
$ws.Range("L2").Value = 82211.88392976951
$ws.Range("M2").Value = 21558.28453827675
$ws.Range("N2").Value = 10979.3631303427
$ws.Range("O2").Value = 9591.909189750015

# --- 2035 ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 1990.187149144482
$ws.Range("B2").Value = 6771.312033453911
$ws.Range("E2").Value = 68594.62611205096
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 49425.72334238004
$ws.Range("L2").Value = 82211.88392976951
$ws.Range("M2").Value = 26927.48989038075
$ws.Range("N2").Value = 16039.77123262105
$ws.Range("O2").Value = 15303.0518371641

# --- 2040 ---
$ws = $wb.Worksheets.Item("2040")
$ws.Range("A2").Value = 1990.187149144482
$ws.Range("B2").Value = 6771.312033453911
$ws.Range("E2").Value = 68594.62611205096
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 49425.72334238004
$ws.Range("L2").Value = 82211.88392976951
$ws.Range("M2").Value = 26927.48989038075
$ws.Range("N2").Value = 16039.77123262105
$ws.Range("O2").Value = 15303.0518371641

# --- 2045 ---
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 1990.187149144482
$ws.Range("B2").Value = 6771.312033453911
$ws.Range("E2").Value = 68594.62611205096
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 49425.72334238004
$ws.Range("L2").Value = 82211.88392976951
$ws.Range("M2").Value = 26927.48989038075
$ws.Range("N2").Value = 16039.77123262105
$ws.Range("O2").Value = 15303.0518371641

# --- 2050 ---
$ws = $wb.Worksheets.Item("2050")
$ws.Range("A2").Value = 1990.187149144482
$ws.Range("B2").Value = 6771.312033453911
$ws.Range("E2").Value = 68594.62611205096
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 49425.72334238004
$ws.Range("L2").Value = 82211.88392976951
$ws.Range("M2").Value = 26927.48989038075
$ws.Range("N2").Value = 16039.77123262105
$ws.Range("O2").Value = 15303.0518371641
